# "Generate Report for Handoff"
#
# A new handoff was generated for 64ae4af2-4a30-4e99-91db-fb8621460b73.md,
# so its "Latest Handoff Date(time)" cells are refreshed on every sheet
# that tracks it:
#   - zh-cn : Latest Handoff Datetime  02:49:13 -> 02:49:43
#   - de-de : Latest Handoff Datetime  02:49:18 -> 02:49:47
#   - Overview : Latest Handoff Date   = max(zh-cn, de-de) -> 02:49:47

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E6").Value = "2016-03-25 02:49:43"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E6").Value = "2016-03-25 02:49:47"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("D6").Value = "2016-03-25 02:49:47"
